$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 130937854
$ws.Range("B3").Value = 57881
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = 'Spillkråka'
$ws.Range("G3").Value = 'Dryocopus martius'
$ws.Range("Q3").Value = 489668
$ws.Range("R3").Value = 7004128
$ws.Range("AC3").Value = 'Rejäla hackspår, färska och äldre, I två levande granar och i ytlig grov rotdel.'

# Row 4
$ws.Range("A4").Value = 130937843
$ws.Range("M4").Value = 'färska spår'
$ws.Range("Q4").Value = 489760
$ws.Range("R4").Value = 7004232
$ws.Range("AC4").Value = 'Ringhack, färska och äldre, i riklig mängd längs flera meter högt upp på en granstam med spår av rikligt sav/kådaflöde.'

# Row 5
$ws.Range("A5").Value = 130937852
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("M5").Value = 'äldre spår'
$ws.Range("Q5").Value = 489520
$ws.Range("R5").Value = 7004161
$ws.Range("AC5").Value = 'Ringhack, äldre, ytliga enstaka längs flera meter på en granstam vid kanten mot yngre skog.'
$ws.Range("AH5").Value = 'Granskog'
$ws.Range("AJ5").Value = 'gran'
$ws.Range("AK5").Value = 'Picea abies'
$ws.Range("AM5").Value = 'Trädstam på levande träd'
$ws.Range("AO5").Value = 'Stem on living tree # Picea abies'
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = ''
$ws.Range("K5").Value = ''

# Row 6
$ws.Range("A6").Value = 130937857
$ws.Range("B6").Value = 97879
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 221945
$ws.Range("F6").Value = 'Revlummer'
$ws.Range("G6").Value = 'Lycopodium annotinum'
$ws.Range("H6").Value = 'L.'
$ws.Range("Q6").Value = 489680
$ws.Range("R6").Value = 7004154
$ws.Range("M6").Value = ''
$ws.Range("AC6").Value = ''
$ws.Range("AJ6").Value = ''
$ws.Range("AK6").Value = ''
$ws.Range("AM6").Value = ''
$ws.Range("AO6").Value = ''

# Row 7
$ws.Range("A7").Value = 130937863
$ws.Range("B7").Value = 99014
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("I7").Value = "'8"
$ws.Range("J7").Value = 'plantor/tuvor'
$ws.Range("K7").Value = 'fullt utvecklade blad'
$ws.Range("Q7").Value = 489799
$ws.Range("R7").Value = 7004245
$ws.Range("AC7").Value = 'Minst 8 plantor inom ca 1 m2 yta. Grävdes varsamt fram under snötäcket. Det finns sannolikt betydligt mer knärot på fyndplatsen och i skogsbeståndet där fyndplatsen ligger.'
$ws.Range("AH7").Value = 'Barrskog'

# Row 19
$ws.Range("A19").Value = 130937846
$ws.Range("B19").Value = 57884
$ws.Range("D19").Value = 'NT'
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = 'Tretåig hackspett'
$ws.Range("G19").Value = 'Picoides tridactylus'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("M19").Value = 'färska spår'
$ws.Range("Q19").Value = 489591
$ws.Range("R19").Value = 7004206
$ws.Range("AC19").Value = 'Ringhack, främst färska och några äldre, längs flera meter på en granstam vid kant mot yngre skog. Fyndplatsen finns i barrblandskog intill några rotvältor.'
$ws.Range("AJ19").Value = 'gran'
$ws.Range("AK19").Value = 'Picea abies'
$ws.Range("AM19").Value = 'Trädstam på levande träd'
$ws.Range("AO19").Value = 'Stem on living tree # Picea abies'

# Row 20
$ws.Range("A20").Value = 130937860
$ws.Range("B20").Value = 97879
$ws.Range("D20").Value = 'LC'
$ws.Range("E20").Value = 221945
$ws.Range("F20").Value = 'Revlummer'
$ws.Range("G20").Value = 'Lycopodium annotinum'
$ws.Range("H20").Value = 'L.'
$ws.Range("Q20").Value = 489614
$ws.Range("R20").Value = 7004216
$ws.Range("AC20").Value = 'Växer här i barrblandskog.'
$ws.Range("M20").Value = ''
$ws.Range("AJ20").Value = ''
$ws.Range("AK20").Value = ''
$ws.Range("AM20").Value = ''
$ws.Range("AO20").Value = ''
